# Daily attendance processing - swap the order of "System" and the
# recorder's email address in the "Recorded By" column (G) wherever both
# show up together as "System, <email>", turning it into "<email>, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    if ($cell.Text -eq $oldText) {
        $cell.Value = $newText
    }
}
